$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.253.38'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.47'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6703'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07447'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2965'
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.78'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07727'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.799.04'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6790'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.51'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.187'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008312'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.066.84'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.09'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.56'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.205'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.99'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.706'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("E26").Value = '  -4.26%  '
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.515'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.091'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.191'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05327'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.894'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7596'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.147'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.677'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.335.94'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01809'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.733'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9261'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.955'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.39'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08131'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +14.08%  '
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5155'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.770'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.926.15'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.210'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05934'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.13%  '
